$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (it used to wrap the "Норма
#    водоспоживання ... Всього до оплати" block).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists('_GoBack')) {
    $d.Bookmarks('_GoBack').Delete()
}

# ---------------------------------------------------------------------------
# 2. In the "${water}  м3, " line, make "м3" bold (splitting the trailing
#    ", " - which keeps its original formatting incl. the non-breaking
#    space - into its own run).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute('${water}')
if (-not $found) { throw 'could not find ${water}' }
$r.Collapse(0)

# Move past the two spaces right after "${water}" onto the "м" character.
$scan = $r.Duplicate
$scan.End = $scan.Start + 40
$segText = $scan.Text
$mIdx = $segText.IndexOf('м')
if ($mIdx -lt 0) { throw 'could not find м after ${water}' }
$digitIdx = $mIdx + 1

$mRange = $d.Range($r.Start + $mIdx, $r.Start + $digitIdx + 1)
Write-Host 'water mRange text:' $mRange.Text
$mRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. In the "${watering}   м3 " line, make "м3" bold, wrap a new "_GoBack"
#    bookmark around it, and keep the trailing " " as its own (non-bold) run.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute('${watering}')
if (-not $found2) { throw 'could not find ${watering}' }
$r2.Collapse(0)

$scan2 = $r2.Duplicate
$scan2.End = $scan2.Start + 40
$segText2 = $scan2.Text
$mIdx2 = $segText2.IndexOf('м')
if ($mIdx2 -lt 0) { throw 'could not find м after ${watering}' }
$digitIdx2 = $mIdx2 + 1

$mRange2 = $d.Range($r2.Start + $mIdx2, $r2.Start + $digitIdx2 + 1)
Write-Host 'watering mRange text:' $mRange2.Text
$mRange2.Font.Bold = $true

$d.Bookmarks.Add('_GoBack', $mRange2)
